# ------------------------------------------------------------------
# "Updated info, updated N22"
#
#   1. Re-highlights two "all active" rows on N23-Ranges (row 14
#      fully, row 15 partially C:J range) using the same highlighted
#      fill already used elsewhere on that sheet.
#   2. Adds a new worksheet "MainBoard-Relais" (after "N23-Ranges")
#      describing the relay <-> VDC-range mapping (200mV..2000V).
#   3. Leaves the selection/active sheet on the new sheet, matching
#      the saved workbook view.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. N23-Ranges: recolor row 14 (C:J) and row 15 (F:J) to the
#    "highlighted" look already used by the other "V" cells on the
#    sheet. PasteSpecial(formats-only) reuses the existing fill
#    style instead of minting a new one.
# ------------------------------------------------------------------
$ranges = $wb.Worksheets.Item("N23-Ranges")

$highlightSample = $ranges.Range("C8")
$highlightSample.Copy()
$ranges.Range("C14:J14").PasteSpecial(-4122)
$ranges.Range("F15:J15").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# 2. Add the new "MainBoard-Relais" worksheet right after
#    "N23-Ranges".
# ------------------------------------------------------------------
$relais = $wb.Worksheets.Add([Type]::Missing, $ranges)
$relais.Name = "MainBoard-Relais"

$relais.Columns.Item(1).ColumnWidth = 5.140625
$relais.Rows.Item(2).RowHeight = 18.75

# Header row (row 3) first -- this is the order the shared strings
# were actually minted in (200mV..2000V before the "VDC" title).
$relais.Range("C3").Value = "200mV"
$relais.Range("D3").Value = "2V"
$relais.Range("E3").Value = "20V"
$relais.Range("F3").Value = "200V"
$relais.Range("G3").Value = "2000V"
$relais.Range("C3:G3").Font.Bold = $true
$relais.Range("C3:G3").HorizontalAlignment = -4108

$relais.Range("B3").Value = "Range"
$relais.Range("B3").Font.Italic = $true
$relais.Range("B3").HorizontalAlignment = -4152

# Title row (merged C2:G2), bold/size-14, centered.
$relais.Range("C2").Value = "VDC"
$title = $relais.Range("C2:G2")
$title.Merge()
$title.Font.Bold = $true
$title.Font.Size = 14
$title.HorizontalAlignment = -4108

$relais.Range("B4").Value = "Relais"
$relais.Range("B4").Font.Italic = $true
$relais.Range("B4").HorizontalAlignment = -4152

# Column A (rows 4-26) + column B (rows 5-28) use the bold font used
# for row markers throughout the sheet.
$relais.Range("A4:A26").Font.Bold = $true
$relais.Range("B5:B28").Font.Bold = $true

# Relay numbers, column B, rows 5-26.
$relayNumbers = @(1202,1203,1204,1301,1302,1303,1304,1305,1306,1307,1401,1402,1601,1602,1603,1604,1605,1606,1607,1608,1609,1610)
$row = 5
foreach ($n in $relayNumbers) {
    $relais.Cells.Item($row, 2).Value = $n
    $row++
}

# Rows that carry per-range "x" / "V" marks (centered). Default to
# "x" everywhere first, then overwrite the active range with "V".
$xRows = @(5,6,7,10,25)
foreach ($r in $xRows) {
    $relais.Range("C$r`:G$r").Value = "x"
}
$relais.Range("C5:G25").HorizontalAlignment = -4108

$relais.Range("E5").Value = "V"
$relais.Range("F6").Value = "V"
$relais.Range("G7").Value = "V"
$relais.Range("C10").Value = "V"
$relais.Range("D10").Value = "V"

$highlightSample.Copy()
$relais.Range("E5").PasteSpecial(-4122)
$relais.Range("F6").PasteSpecial(-4122)
$relais.Range("G7").PasteSpecial(-4122)
$relais.Range("C10:D10").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Rows 8, 9 and 13 keep the centered (but empty / unfilled) look.
$relais.Range("C8:G9").HorizontalAlignment = -4108
$relais.Range("C13:G13").HorizontalAlignment = -4108

$relais.Range("D13").Select()
$relais.Activate()
